$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2 through 43
# from 2025-06-17 (45825) to 2025-06-18 (45826)
$ws.Range("C2:C43").Value = 45826
